# Update the build/version timestamp embedded in the "About" sheet and in
# the "build_version" column of the "Boundaries and methane sources" sheet.
#
# Old: mines - January 30 (built on January 30 2026 16.19.47 EST)
# New: mines - January 30 (built on February 02 2026 12.49.33 EST)

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

# --- "About" sheet -------------------------------------------------------
$aboutWs = $wb.Worksheets.Item("About")

$a2 = $aboutWs.Range("A2").Value()
$aboutWs.Range("A2").Value = $a2.Replace($oldStamp, $newStamp)

$a6 = $aboutWs.Range("A6").Value()
$aboutWs.Range("A6").Value = $a6.Replace($oldStamp, $newStamp)

# --- "Boundaries and methane sources" sheet -------------------------------
$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")

$lastRow = $dataWs.Cells.Item(1, 1).End(4).Row
if ($lastRow -lt 14) { $lastRow = 14 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $dataWs.Cells.Item($r, 19)  # column S = build_version
    $v = $cell.Value()
    if ($v -ne $null -and $v -ne "") {
        $cell.Value = $v.Replace($oldStamp, $newStamp)
    }
}
